# New weekly data point: insert two rows at the top of the data block
# (rows 3 and 4), pushing the existing records down. This mirrors the
# "Fruta / hortaliza, semanal" update where a new week's Ciboulette
# prices (Primera / Segunda) are prepended to the series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("3:4").Insert()

# Row 3: Primera calidad, latest week
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = 45149
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 100112039
$ws.Range("G3").Value = "Ciboulette"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 2500
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = 2500
$ws.Range("N3").Value = "`$/docena de atados"
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 833
$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = "Hortaliza"

# Row 4: Segunda calidad, latest week
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 45149
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100112039
$ws.Range("G4").Value = "Ciboulette"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Segunda"
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = 2000
$ws.Range("N4").Value = "`$/docena de atados"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 667
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = "Hortaliza"
